# The commit swaps the two theme parts of the deck: the theme that used to
# live at ppt/theme/theme1.xml ("Office Theme" / "Office" colour scheme) and
# the theme that used to live at ppt/theme/theme2.xml ("Integral" / "Red
# Violet" colour scheme) trade their colour palettes - theme2.xml (the
# design actually bound to the slide master / presentation) ends up with
# the old "Office" colours, while the notes-master-only theme1.xml ends up
# with the old "Red Violet" colours. Font scheme and format scheme are
# byte-for-byte identical between the two parts, so the only observable
# difference is the 10 colour slots that differ (dk1/lt1 - black/white -
# are shared by both palettes already).

function ConvertTo-BgrColor([string]$hex) {
    # PowerPoint's RGB COM property stores colours as 0x00BBGGRR.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p = $ppt.ActivePresentation

# Target palette: the colours that used to belong to the "Office Theme" /
# "Office" colour scheme (theme1.xml before the edit), in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

# The presentation's live theme (bound to the slide master / all slides) is
# reached through any slide's ThemeColorScheme, which exposes all twelve
# colour-scheme slots (unlike Master.ColorScheme, which only exposes eight
# and clobbers the scheme's name when written to).
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-BgrColor $officeColors[$i - 1]
}
